$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 135.53
$ws.Range("I15").Value = 135.53
$ws.Range("K15").Value = 406.59
$ws.Range("M15").Value = -237.59
$ws.Range("H123").Value = 29075.428
$ws.Range("J123").Value = 29075.428
$ws.Range("L123").Value = 29075.428
$ws.Range("N123").Value = -38875.428
$ws.Range("H129").Value = 2200.682
$ws.Range("J129").Value = 1932.5264
$ws.Range("L129").Value = 5797.5792
$ws.Range("N129").Value = -15797.5792

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 565.8461
$ws.Range("I5").Value = 541.125
$ws.Range("J5").Value = 605.4
$ws.Range("K5").Value = 541.125
$ws.Range("L5").Value = 605.4
$ws.Range("M5").Value = -429.125
$ws.Range("N5").Value = -829.4
$ws.Range("H61").Value = 2643.3877
$ws.Range("I61").Value = 1230.5217
$ws.Range("J61").Value = 3893.2307
$ws.Range("K61").Value = 1230.5217
$ws.Range("L61").Value = 3893.2307
$ws.Range("M61").Value = -1018.5217
$ws.Range("N61").Value = -4317.2307
$ws.Range("H74").Value = 3306.5
$ws.Range("I74").Value = 3720.2354
$ws.Range("J74").Value = 1899.8
$ws.Range("K74").Value = 3720.2354
$ws.Range("L74").Value = 1899.8
$ws.Range("M74").Value = -2846.2354
$ws.Range("N74").Value = -3647.8
$ws.Range("H77").Value = 3306.5
$ws.Range("I77").Value = 3720.2354
$ws.Range("J77").Value = 1899.8
$ws.Range("K77").Value = 18601.177
$ws.Range("L77").Value = 9499
$ws.Range("M77").Value = -14233.177
$ws.Range("N77").Value = -18235
$ws.Range("H131").Value = 50418.75
$ws.Range("J131").Value = 50418.75
$ws.Range("L131").Value = 50418.75
$ws.Range("N131").Value = -60498.75
$ws.Range("H136").Value = 2643.3877
$ws.Range("I136").Value = 1230.5217
$ws.Range("J136").Value = 3893.2307
$ws.Range("K136").Value = 3691.5651
$ws.Range("L136").Value = 11679.6921
$ws.Range("M136").Value = -1141.5651
$ws.Range("N136").Value = -16779.6921

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 565.8461
$ws.Range("I4").Value = 541.125
$ws.Range("J4").Value = 605.4
$ws.Range("K4").Value = 541.125
$ws.Range("L4").Value = 605.4
$ws.Range("M4").Value = -426.125
$ws.Range("N4").Value = -835.4
$ws.Range("H44").Value = 50000
$ws.Range("J44").Value = 50000
$ws.Range("L44").Value = 50000
$ws.Range("N44").Value = -50994
$ws.Range("H107").Value = 2006.8387
$ws.Range("I107").Value = 1899.9642
$ws.Range("J107").Value = 3004.3333
$ws.Range("K107").Value = 1899.9642
$ws.Range("L107").Value = 3004.3333
$ws.Range("M107").Value = 20.03580000000011
$ws.Range("N107").Value = -6844.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4628.486
$ws.Range("I31").Value = 1612.0834
$ws.Range("J31").Value = 6202.2607
$ws.Range("K31").Value = 1612.0834
$ws.Range("L31").Value = 6202.2607
$ws.Range("M31").Value = -1317.0834
$ws.Range("N31").Value = -6792.2607
$ws.Range("H34").Value = 4628.486
$ws.Range("I34").Value = 1612.0834
$ws.Range("J34").Value = 6202.2607
$ws.Range("K34").Value = 1612.0834
$ws.Range("L34").Value = 6202.2607
$ws.Range("M34").Value = -1410.0834
$ws.Range("N34").Value = -6606.2607
$ws.Range("H41").Value = 24865.572
$ws.Range("I41").Value = 5000
$ws.Range("J41").Value = 32811.8
$ws.Range("K41").Value = 5000
$ws.Range("L41").Value = 32811.8
$ws.Range("M41").Value = -4572
$ws.Range("N41").Value = -33667.8
$ws.Range("H50").Value = 40675
$ws.Range("J50").Value = 40675
$ws.Range("L50").Value = 40675
$ws.Range("N50").Value = -41925
$ws.Range("H51").Value = 250019900
$ws.Range("J51").Value = 39800
$ws.Range("L51").Value = 39800
$ws.Range("N51").Value = -41272
$ws.Range("H59").Value = 35118.57
$ws.Range("J59").Value = 35118.57
$ws.Range("L59").Value = 35118.57
$ws.Range("N59").Value = -37408.57
$ws.Range("H60").Value = 21288.955
$ws.Range("I60").Value = 1250
$ws.Range("J60").Value = 23292.85
$ws.Range("K60").Value = 1250
$ws.Range("L60").Value = 23292.85
$ws.Range("M60").Value = -739
$ws.Range("N60").Value = -24314.85
$ws.Range("H61").Value = 250019900
$ws.Range("J61").Value = 39800
$ws.Range("L61").Value = 39800
$ws.Range("N61").Value = -40496
$ws.Range("H68").Value = 57000
$ws.Range("J68").Value = 57000
$ws.Range("L68").Value = 57000
$ws.Range("N68").Value = -58498
$ws.Range("H71").Value = 57000
$ws.Range("J71").Value = 57000
$ws.Range("L71").Value = 171000
$ws.Range("N71").Value = -178488
$ws.Range("H74").Value = 13000
$ws.Range("J74").Value = 13000
$ws.Range("L74").Value = 13000
$ws.Range("N74").Value = -14748
$ws.Range("H77").Value = 13000
$ws.Range("J77").Value = 13000
$ws.Range("L77").Value = 39000
$ws.Range("N77").Value = -47736
$ws.Range("H132").Value = 34370.188
$ws.Range("I132").Value = 1438.1177
$ws.Range("K132").Value = 4314.3531
$ws.Range("M132").Value = -1784.3531
$ws.Range("H134").Value = 1637.6296
$ws.Range("I134").Value = 928.85
$ws.Range("J134").Value = 3662.7144
$ws.Range("K134").Value = 2786.55
$ws.Range("L134").Value = 10988.1432
$ws.Range("M134").Value = -251.5500000000002
$ws.Range("N134").Value = -16058.1432

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 40008084
$ws.Range("I137").Value = 2698.4614
$ws.Range("J137").Value = 83347256
$ws.Range("K137").Value = 8095.3842
$ws.Range("L137").Value = 250041768
$ws.Range("M137").Value = -2995.3842
$ws.Range("N137").Value = -250051968

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H38").Value = 14947.5
$ws.Range("I38").Value = 14800
$ws.Range("J38").Value = 14996.667
$ws.Range("K38").Value = 14800
$ws.Range("L38").Value = 14996.667
$ws.Range("M38").Value = -14337
$ws.Range("N38").Value = -15922.667
$ws.Range("H43").Value = 18097
$ws.Range("I43").Value = 15250
$ws.Range("J43").Value = 19995
$ws.Range("K43").Value = 15250
$ws.Range("L43").Value = 19995
$ws.Range("M43").Value = -15099
$ws.Range("N43").Value = -20297
$ws.Range("H113").Value = 6959.684
$ws.Range("I113").Value = 8261.4
$ws.Range("K113").Value = 8261.4
$ws.Range("M113").Value = -6091.4
$ws.Range("H123").Value = 15287.111
$ws.Range("J123").Value = 15287.111
$ws.Range("L123").Value = 15287.111
$ws.Range("N123").Value = -20187.111

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H94").Value = 43552.445
$ws.Range("J94").Value = 43552.445
$ws.Range("L94").Value = 43552.445
$ws.Range("N94").Value = -44904.445
$ws.Range("H98").Value = 39398.4
$ws.Range("J98").Value = 39398.4
$ws.Range("L98").Value = 39398.4
$ws.Range("N98").Value = -45388.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H42").Value = 22000
$ws.Range("J42").Value = 22000
$ws.Range("L42").Value = 22000
$ws.Range("N42").Value = -22756
$ws.Range("H43").Value = 19150
$ws.Range("I43").Value = 2900
$ws.Range("J43").Value = 29983.334
$ws.Range("K43").Value = 2900
$ws.Range("L43").Value = 29983.334
$ws.Range("M43").Value = -2751
$ws.Range("N43").Value = -30281.334
$ws.Range("H104").Value = 47376.8
$ws.Range("J104").Value = 47376.8
$ws.Range("L104").Value = 47376.8
$ws.Range("N104").Value = -54364.8
$ws.Range("H132").Value = 2758.3635
$ws.Range("I132").Value = 2046.7646
$ws.Range("J132").Value = 3514.4375
$ws.Range("K132").Value = 6140.293799999999
$ws.Range("L132").Value = 10543.3125
$ws.Range("M132").Value = -3610.293799999999
$ws.Range("N132").Value = -15603.3125
$ws.Range("H136").Value = 22024.818
$ws.Range("I136").Value = 74314.336
$ws.Range("J136").Value = 2416.25
$ws.Range("K136").Value = 222943.008
$ws.Range("L136").Value = 7248.75
$ws.Range("M136").Value = -220393.008
$ws.Range("N136").Value = -12348.75
